$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting existing rows 79..117 down to 80..118.
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with the new observation.
$ws.Range("A79").Value = 8
$ws.Range("B79").Value = "Terminal La Palmera de La Serena"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44813
$ws.Range("E79").Value = 4
$ws.Range("F79").Value = 100112052
$ws.Range("G79").Value = "Albahaca"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 1300
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = 4250
$ws.Range("N79").Value = "$/paquete"
$ws.Range("O79").Value = "Región de Arica y Parinacota"
$ws.Range("P79").Value = 4250
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"
